# Add Sergio Cerdeira to the P2P Members list.
# A new row is inserted right before the last existing row (Yasmina Shah
# Esmaeili), which pushes her row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 48, shifting the current row 48 down to 49.
$ws.Rows.Item(48).Insert()

# Populate the new member's row: MEMBERS | Country | Affiliation | Topic | email
$ws.Range("A48").Value = "[Sergio Cerdeira](https://oceanexpert.org/expert/30635)"
$ws.Range("B48").Value = "Mexico"
$ws.Range("C48").Value = "[CONABIO - SIMAR](https://simar.conabio.gob.mx/)"
$ws.Range("D48").Value = "RS"
$ws.Range("E48").Value = "scerdeira@conabio.gob.mx"

$ws.Range("C51").Select()
